$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking identifier strings must be forced to Text
# format before assignment, otherwise Excel auto-converts them to numbers.
$idCells = @("B2", "K2", "M2", "O2", "Q2", "B3", "K3", "M3", "O3", "Q3", "B4", "E4", "H4", "K4", "M4", "O4", "Q4")
foreach ($addr in $idCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# ---- Row 2 ----
$ws.Range("A2").Value = 48
$ws.Range("B2").Value = "3244307"
$ws.Range("C2").Value = "tiotropium bromide"
$ws.Range("D2").Value = "Chemicals & Drugs"
$ws.Range("E2").Value = "4027653"
$ws.Range("F2").Value = "patients"
$ws.Range("G2").Value = "Living Beings"
$ws.Range("H2").Value = "3763631"
$ws.Range("I2").Value = "bullous eruption"
$ws.Range("J2").Value = "Disorders"
$ws.Range("K2").Value = "119763642"
$ws.Range("L2").Value = "is not administered to"
$ws.Range("M2").Value = "58220902"
$ws.Range("N2").Value = "is administered to"
$ws.Range("O2").Value = "53071987"
$ws.Range("P2").Value = "treats"
$ws.Range("Q2").Value = "51266807"
$ws.Range("R2").Value = "is process of"

# ---- Row 3 ----
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "3244307"
$ws.Range("C3").Value = "tiotropium bromide"
$ws.Range("D3").Value = "Chemicals & Drugs"
$ws.Range("E3").Value = "4033983"
$ws.Range("F3").Value = "pharmaceutical preparations"
$ws.Range("G3").Value = "Chemicals & Drugs"
$ws.Range("H3").Value = "3763631"
$ws.Range("I3").Value = "bullous eruption"
$ws.Range("J3").Value = "Disorders"
$ws.Range("K3").Value = "89903399"
$ws.Range("L3").Value = "coexists with"
$ws.Range("M3").Value = "64510312"
$ws.Range("N3").Value = "is a"
$ws.Range("O3").Value = "116913534"
$ws.Range("P3").Value = "interacts with"
$ws.Range("Q3").Value = "76350720"
$ws.Range("R3").Value = "causes"

# ---- Row 4 ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "3244307"
$ws.Range("C4").Value = "tiotropium bromide"
$ws.Range("D4").Value = "Chemicals & Drugs"
$ws.Range("E4").Value = "5648228"
$ws.Range("F4").Value = "therapeutic procedure"
$ws.Range("G4").Value = "Procedures"
$ws.Range("H4").Value = "3763631"
$ws.Range("I4").Value = "bullous eruption"
$ws.Range("J4").Value = "Disorders"
$ws.Range("K4").Value = "61122555"
$ws.Range("L4").Value = "uses"
$ws.Range("M4").Value = "55053777"
$ws.Range("N4").Value = "is compared with"
$ws.Range("O4").Value = "128974761"
$ws.Range("P4").Value = "is compared with"
$ws.Range("Q4").Value = "66163822"
$ws.Range("R4").Value = "treats"
